$wb = $excel.ActiveWorkbook

$wsWeek6 = $wb.Worksheets.Item("Week 6")
$wsWeek7 = $wb.Worksheets.Item("Week 7")

# Copy the date/time/description number formats from an existing filled
# row (Week 6, row 4) onto the new rows being added in Week 7 (rows 2-4),
# matching the formatting pattern already used throughout the workbook.
$wsWeek6.Range("A4:D4").Copy()
$wsWeek7.Range("A2:D2").PasteSpecial(-4122)
$wsWeek7.Range("A3:D3").PasteSpecial(-4122)
$wsWeek7.Range("A4:D4").PasteSpecial(-4122)

# Enter the two new timesheet entries for Week 7.
$wsWeek7.Range("A2").Value2 = 43516
$wsWeek7.Range("B2").Value2 = 0.52083333333333337
$wsWeek7.Range("C2").Value2 = 0.63541666666666663
$wsWeek7.Range("D2").Value2 = "Worked on Admin Area Content"
$wsWeek7.Range("E2").Value2 = 2.75

$wsWeek7.Range("A3").Value2 = 43517
$wsWeek7.Range("B3").Value2 = 0.5
$wsWeek7.Range("C3").Value2 = 0.54166666666666663
$wsWeek7.Range("D3").Value2 = "Worked on Admin Area Content"
$wsWeek7.Range("E3").Value2 = 1

# Select D5 on Week 7 and make it the active sheet/tab, matching the
# author's last recorded position when they saved the workbook.
$wsWeek7.Activate()
$wsWeek7.Range("D5").Select()
